$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "agentNum" (column D) values for the affected row blocks ---
# Block 1: rows 2-36  -> 2
$ws.Range("D2:D36").Value = 2
# Block 2: rows 37-71 -> 1
$ws.Range("D37:D71").Value = 1
# Block 3: rows 72-112 -> 1
$ws.Range("D72:D112").Value = 1
# Block 4: rows 113-153 stay at 2 (no change)
# Block 5: rows 154-194 -> 3
$ws.Range("D154:D194").Value = 3
# Block 6: rows 195-235 -> 4
$ws.Range("D195:D235").Value = 4
# Block 7: rows 236-276 -> 1
$ws.Range("D236:D276").Value = 1
# Block 8: rows 277-317 -> 4
$ws.Range("D277:D317").Value = 4
# Block 9: rows 318-358 -> 3
$ws.Range("D318:D358").Value = 3
# Block 10: rows 359-399 stay at 2 (no change)

# --- Apply an AutoFilter on column B (study) showing only values 2 and 3 ---
$ws.Range("A1:S399").AutoFilter(2, @("2","3"), 7)

# --- Restore the view: scroll position + selection ---
$ws.Range("D277:D399").Select()
$excel.ActiveWindow.ScrollRow = 78
